# The sheet originally has columns:
#   A name | B address | C website | D phone_number | E reviews_count |
#   F reviews_average | G latitude | H longitude | I is_permanently_closed |
#   J gmaps_link | K latest_review_date
#
# The "reviews_count" column (E) is being removed entirely, shifting every
# column to its right (F:K) one place to the left (-> E:J). Deleting the
# whole column reproduces exactly that: remaining data slides left and the
# used range shrinks from A1:K21 to A1:J21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
